# Updates the cryptos price/volume snapshot (daily GitHub Actions refresh),
# plus a ShibaInu/ImmutableX rank swap at rows 21-22.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.700.98"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "2.506.57"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.04%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "322.43"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "107.90"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.03%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.560"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.96%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "40.22"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("E11").Value = "  -0.38%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "19.51"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.17"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "2.899.51"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "2.510.26"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "47.621.08"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.34"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.24%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "2.77"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.16%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.94"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "247.19"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  -0.01%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "25.75"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.23"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("E29").Value = "  +4.13%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "34.80"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -9.42%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "49.82"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "20.00"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -1.61%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0782"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("E37").Value = "  -1.23%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "4.68"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("E40").Value = "  -0.30%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "22.24"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.78%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.30%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "118.81"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "2.002.21"
$ws.Range("E45").Value = "  +0.60%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("E48").Value = "  +0.52%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.08"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -3.43%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "56.58"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
